$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (interested count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 7736
$wsExhibit.Range("F5").Value = 5630
$wsExhibit.Range("F7").Value = 75
$wsExhibit.Range("F11").Value = 230

# Sheet "全部类型" (all types) - same events appear here with different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7736
$wsAll.Range("F5").Value = 5630
$wsAll.Range("F7").Value = 75
$wsAll.Range("F13").Value = 230
